$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns ("Staff" and "Admin") ahead of the existing Budget column,
# pushing Budget -> G and Rank -> H.
$ws.Range("E1:F1").EntireColumn.Insert()

# Header row: new column titles + rename old "Rank" header to "Russell"
$ws.Range("E1").Value = "Staff"
$ws.Range("F1").Value = "Admin"
$ws.Range("H1").Value = "Russell"

# Populate new Staff (E) / Admin (F) figures and the new Russell Group flag (H)
# for every data row. Cells with no source value are intentionally left blank.
$ws.Range("E2").Value = 2942
$ws.Range("F2").Value = 4003
$ws.Range("H2").Value = 1
$ws.Range("E3").Value = 4589
$ws.Range("F3").Value = 6107
$ws.Range("H3").Value = 1
$ws.Range("E4").Value = 1137
$ws.Range("F4").Value = 1576
$ws.Range("H4").Value = 0
$ws.Range("E5").Value = 1086
$ws.Range("F5").Value = 1489
$ws.Range("H5").Value = 0
$ws.Range("F6").Value = 3200
$ws.Range("H6").Value = 0
$ws.Range("E7").Value = 1725
$ws.Range("F7").Value = 2515
$ws.Range("H7").Value = 1
$ws.Range("E8").Value = 7700
$ws.Range("F8").Value = 5375
$ws.Range("H8").Value = 1
$ws.Range("E9").Value = 7913
$ws.Range("F9").Value = 3615
$ws.Range("H9").Value = 1
$ws.Range("E10").Value = 7000
$ws.Range("H10").Value = 1
$ws.Range("E11").Value = 2610
$ws.Range("F11").Value = 4033
$ws.Range("H11").Value = 1
$ws.Range("E12").Value = 4390
$ws.Range("F12").Value = 4075
$ws.Range("H12").Value = 1
$ws.Range("E13").Value = 5220
$ws.Range("F13").Value = 3485
$ws.Range("H13").Value = 1
$ws.Range("E14").Value = 3849
$ws.Range("H14").Value = 1
$ws.Range("E15").Value = 3285
$ws.Range("F15").Value = 6199
$ws.Range("H15").Value = 1
$ws.Range("E16").Value = 4020
$ws.Range("H16").Value = 1
$ws.Range("E17").Value = 3235
$ws.Range("F17").Value = 4620
$ws.Range("H17").Value = 1
$ws.Range("E18").Value = 1935
$ws.Range("F18").Value = 3091
$ws.Range("H18").Value = 1
$ws.Range("E19").Value = 3495
$ws.Range("H19").Value = 1
$ws.Range("E20").Value = 1410
$ws.Range("F20").Value = 1805
$ws.Range("H20").Value = 0
$ws.Range("E21").Value = 3330
$ws.Range("F21").Value = 5739
$ws.Range("H21").Value = 1
$ws.Range("F22").Value = 1872
$ws.Range("H22").Value = 0
$ws.Range("E23").Value = 2414
$ws.Range("F23").Value = 1489
$ws.Range("H23").Value = 1
$ws.Range("F24").Value = 3290
$ws.Range("H24").Value = 0

# Column widths (closest representable values under this engines 1/6-character
# width quantization) matching the post-edit layout.
$ws.Columns.Item(5).ColumnWidth = 4.593
$ws.Columns.Item(6).ColumnWidth = 5.76
$ws.Columns.Item(7).ColumnWidth = 6.26
$ws.Columns.Item(8).ColumnWidth = 6.427

# Restore the selected cell shown in the saved view
$null = $ws.Range("E43").Select()
